$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 81, pushing existing rows 81-94 down to 83-96.
$ws.Rows.Item(81).Insert()
$ws.Rows.Item(81).Insert()

# Populate new row 81 (Hass / Segunda, Perú origin)
$ws.Cells.Item(81, 1).Value = 1
$ws.Cells.Item(81, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(81, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(81, 4).Value = 44711
$ws.Cells.Item(81, 5).Value = 15
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100106
$ws.Cells.Item(81, 8).Value = "Oleaginosos"
$ws.Cells.Item(81, 9).Value = 100106002
$ws.Cells.Item(81, 10).Value = "Palta"
$ws.Cells.Item(81, 11).Value = "Hass"
$ws.Cells.Item(81, 12).Value = "Segunda"
$ws.Cells.Item(81, 13).Value = 520
$ws.Cells.Item(81, 14).Value = 15000
$ws.Cells.Item(81, 15).Value = 16000
$ws.Cells.Item(81, 16).Value = 15500
$ws.Cells.Item(81, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(81, 18).Value = "Perú"
$ws.Cells.Item(81, 19).Value = 1550
$ws.Cells.Item(81, 20).Value = 10

# Populate new row 82 (Hass / Tercera, Perú origin)
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value = 44711
$ws.Cells.Item(82, 5).Value = 15
$ws.Cells.Item(82, 6).Value = "Fruta"
$ws.Cells.Item(82, 7).Value = 100106
$ws.Cells.Item(82, 8).Value = "Oleaginosos"
$ws.Cells.Item(82, 9).Value = 100106002
$ws.Cells.Item(82, 10).Value = "Palta"
$ws.Cells.Item(82, 11).Value = "Hass"
$ws.Cells.Item(82, 12).Value = "Tercera"
$ws.Cells.Item(82, 13).Value = 520
$ws.Cells.Item(82, 14).Value = 15000
$ws.Cells.Item(82, 15).Value = 16000
$ws.Cells.Item(82, 16).Value = 15500
$ws.Cells.Item(82, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(82, 18).Value = "Perú"
$ws.Cells.Item(82, 19).Value = 1550
$ws.Cells.Item(82, 20).Value = 10
